# Mold Man app, v0.2 commit
# Builds the "mold management" lookup sheet: a set of filter/label cells
# (rows 2-5) plus a data-table header row (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Cell text content.
# Order matters: it reproduces the shared-string table order of the
# authored workbook (header row typed first, then the filter labels
# column-by-column, then the "~" range separators last).
# ---------------------------------------------------------------------

# Row 10: data table header
$ws.Range("A10").Value = "管理籍"
$ws.Range("B10").Value = "製品名"
$ws.Range("C10").Value = "中型"
$ws.Range("D10").Value = "外寸長"
$ws.Range("E10").Value = "外寸幅"
$ws.Range("F10").Value = "外寸高"
$ws.Range("G10").Value = "内寸長"
$ws.Range("H10").Value = "内寸幅"
$ws.Range("I10").Value = "内寸深"
$ws.Range("J10").Value = "蓋"
$ws.Range("K10").Value = "蓋付"
$ws.Range("L10").Value = "製造年月"
$ws.Range("M10").Value = "用途"

# Column A filter labels (rows 2-5)
$ws.Range("A2").Value = "管理籍　："
$ws.Range("A3").Value = "外寸（長さ）："
$ws.Range("A4").Value = "外寸（幅）："
$ws.Range("A5").Value = "外寸（高さ）："

# Column F filter labels (rows 3-5)
$ws.Range("F3").Value = "内寸（長さ）："
$ws.Range("F4").Value = "内寸（幅）："
$ws.Range("F5").Value = "内寸（高さ）："

# K3 filter label
$ws.Range("K3").Value = "製品名："

# "~" range separators (columns D and I)
$ws.Range("D3").Value = "～"
$ws.Range("D4").Value = "～"
$ws.Range("D5").Value = "～"
$ws.Range("I3").Value = "～"
$ws.Range("I4").Value = "～"
$ws.Range("I5").Value = "～"

# ---------------------------------------------------------------------
# Formatting.
# Both the header row and the filter-label boxes use the same light
# gray fill ("White, Background 1, Darker 15%" == RGB D9D9D9); the
# header additionally gets a thin box border around every cell.
#
# Each format is first built on a scratch cell, then propagated with
# Copy/PasteSpecial(xlPasteFormats) so every destination cell picks up
# the combined format in a single paste instead of accumulating it
# property-by-property (which would otherwise register an extra,
# intermediate cell style for every single property write).
# ---------------------------------------------------------------------
$fillColor = 14277081   # RGB(217,217,217) packed as BGR for COM .Color
$xlPasteFormats = -4122
$scratch = $ws.Range("Z1")

# Header style: fill + thin border on all four sides of every cell.
$scratch.Interior.Color = $fillColor
$scratch.Borders.LineStyle = 1
$scratch.Borders.Weight = 2
$scratch.Copy()
$ws.Range("A10:M10").PasteSpecial($xlPasteFormats)
$scratch.Clear()

# Label style: fill only, no border.
$scratch.Interior.Color = $fillColor
$scratch.Copy()
$ws.Range("A2:B2").PasteSpecial($xlPasteFormats)
$ws.Range("A3:B3").PasteSpecial($xlPasteFormats)
$ws.Range("F3:G3").PasteSpecial($xlPasteFormats)
$ws.Range("K3:L3").PasteSpecial($xlPasteFormats)
$ws.Range("A4:B4").PasteSpecial($xlPasteFormats)
$ws.Range("F4:G4").PasteSpecial($xlPasteFormats)
$ws.Range("A5:B5").PasteSpecial($xlPasteFormats)
$ws.Range("F5:G5").PasteSpecial($xlPasteFormats)
$scratch.Clear()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# View state: selection moves to I4.
# ---------------------------------------------------------------------
$ws.Range("I4").Select() | Out-Null
